# Minor edits to the color sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "table-1" to "What Color"
$ws.Name = "What Color"

# Update the active selection from N12 to N24
$ws.Range("N24").Select()
